$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (even_MAG-GUT84132.fa) entirely; this shifts row 7
# (even_MAG-GUT84562.fa) up to become the new row 6, matching the diff.
$ws.Rows.Item(6).Delete()
